# Add the "Слайд 17 (Проблемы)" conclusion section to the end of the
# presentation narration document: one Heading1 paragraph followed by
# two body paragraphs, inserted just before the trailing empty
# paragraph that closes the document body.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range

# Make room for the three new paragraphs right before the final
# (empty) paragraph, preserving that trailing paragraph as-is.
$lastRange.InsertParagraphBefore()
$lastRange.InsertParagraphBefore()
$lastRange.InsertParagraphBefore()

$baseIndex = $d.Paragraphs.Count - 3

# --- Heading paragraph: "Слайд 17 (Проблемы)" ---------------------------
$headingPara = $d.Paragraphs.Item($baseIndex)
$headingPara.Range.InsertAfter("Слайд 17 (Проблемы)")
$headingPara.Style = "Heading 1"
$headingPara.Format.KeepWithNext = $false
$headingPara.Format.KeepTogether = $false
$headingPara.Format.SpaceBefore = 24
$headingPara.Format.LineSpacingRule = 5
$headingPara.Range.Font.Bold = $true
$headingPara.Range.Font.Size = 23
$headingPara.Range.Font.SizeBi = 23

$bmRange = $d.Range($headingPara.Range.Start, $headingPara.Range.Start)
$d.Bookmarks.Add("_xaji0y6dpbhs", $bmRange)

# --- Body paragraph 1 -----------------------------------------------------
$body1 = $d.Paragraphs.Item($baseIndex + 1)
$body1.Range.InsertAfter("Основная проблема при тестировании системы это корректность аппроксимации переходного интервала прямой линией. Обратите внимание на первый график на слайде. Если на вход системе поступит маленькое значение k, значения функции разладки на переходном интервале будут меньшезначений аппроксимации на первых k точках после момента Q. Это частично нивелируется уменьшением разности параметров T и L. Однако мне встретились параметры сигнала, на которых даже в случае отсутствия шума вероятность ложноотрицательного обнаружения разладки была равна 1. Функция разладки данной ситуации привидена на втором графике, где частота периодики равна 1/7. На графике видно, что функция неоднородности начинает быстро расти с некоторой задержкой и это поведение нуждается в дальнейшем исследовании.")
$body1.Format.SpaceBefore = 12
$body1.Format.SpaceAfter = 12
$body1.Format.LineSpacingRule = 5

# --- Body paragraph 2 -----------------------------------------------------
$body2 = $d.Paragraphs.Item($baseIndex + 2)
$body2.Range.InsertAfter("Однако, при достаточно большом значении k, например 30, при стандартном отклонении шума, не превышающем половину амплитуды сигнала, вероятность точного обнаружения для частот от 1/3 до 1/9 составляла не менее 98 процентов. А при k равном 45 - эта вероятность составляла 1.")
$body2.Format.SpaceBefore = 12
$body2.Format.SpaceAfter = 12
$body2.Format.LineSpacingRule = 5
